# Update the "RW [b.rw] RT [b.rt]" label on Sheet1 to just "RW [b.rw]"
# and move the active selection to C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("C3").Value = "RW [b.rw]"

$ws.Range("C23").Select()
